$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 96775736
$ws.Range("B2").Value = 89412
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 5442
$ws.Range("F2").Value = 'Tallticka'
$ws.Range("G2").Value = 'Porodaedalea pini'
$ws.Range("H2").Value = '(Brot.) Murrill'
$ws.Range("Q2").Value = 658730.5226168972
$ws.Range("R2").Value = 6637449.43415721

# Row 3
$ws.Range("A3").Value = 96777744
$ws.Range("B3").Value = 98520
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = 'Blåsippa'
$ws.Range("G3").Value = 'Hepatica nobilis'
$ws.Range("H3").Value = 'Schreb.'
$ws.Range("Q3").Value = 658769.3765012868
$ws.Range("R3").Value = 6637283.535847809

# Row 4
$ws.Range("A4").Value = 96778360
$ws.Range("B4").Value = 98520
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = 'Blåsippa'
$ws.Range("G4").Value = 'Hepatica nobilis'
$ws.Range("H4").Value = 'Schreb.'
$ws.Range("Q4").Value = 658743.2312543363
$ws.Range("R4").Value = 6637305.564015599

# Row 5
$ws.Range("A5").Value = 96777066
$ws.Range("B5").Value = 89392
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 1202
$ws.Range("F5").Value = 'Ullticka'
$ws.Range("G5").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H5").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q5").Value = 658712.8441804391
$ws.Range("R5").Value = 6637414.964914286

# Row 6
$ws.Range("A6").Value = 96780597
$ws.Range("B6").Value = 89392
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = 'Ullticka'
$ws.Range("G6").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H6").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q6").Value = 658693.8994370478
$ws.Range("R6").Value = 6637469.000542388

# Row 7
$ws.Range("A7").Value = 96780135
$ws.Range("B7").Value = 89392
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = 'Ullticka'
$ws.Range("G7").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H7").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q7").Value = 658733.0048414731
$ws.Range("R7").Value = 6636920.163120084

# Row 8
$ws.Range("A8").Value = 96779183
$ws.Range("B8").Value = 98520
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 222498
$ws.Range("F8").Value = 'Blåsippa'
$ws.Range("G8").Value = 'Hepatica nobilis'
$ws.Range("H8").Value = 'Schreb.'
$ws.Range("Q8").Value = 658704.2164550385
$ws.Range("R8").Value = 6637062.857129692

# Row 9
$ws.Range("A9").Value = 96780122
$ws.Range("B9").Value = 89832
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 1209
$ws.Range("F9").Value = 'Rynkskinn'
$ws.Range("G9").Value = 'Phlebia centrifuga'
$ws.Range("H9").Value = 'P.Karst.'
$ws.Range("Q9").Value = 658733.0048414731
$ws.Range("R9").Value = 6636920.163120084

# Row 10
$ws.Range("A10").Value = 96779825
$ws.Range("B10").Value = 89376
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 4660
$ws.Range("F10").Value = 'Rävticka'
$ws.Range("G10").Value = 'Inocutis rheades'
$ws.Range("H10").Value = '(Pers.) Fiasson & Niemelä'
$ws.Range("Q10").Value = 658637.0597997338
$ws.Range("R10").Value = 6636982.990721731

# Row 11
$ws.Range("A11").Value = 96779798
$ws.Range("B11").Value = 43464
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 101735
$ws.Range("F11").Value = 'Jättesvampmal'
$ws.Range("G11").Value = 'Scardia boletella'
$ws.Range("H11").Value = '(Fabricius, 1794)'
$ws.Range("Q11").Value = 658637.0597997338
$ws.Range("R11").Value = 6636982.990721731

# Row 12
$ws.Range("A12").Value = 96780278
$ws.Range("B12").Value = 93132
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 2671
$ws.Range("F12").Value = 'Fällmossa'
$ws.Range("G12").Value = 'Antitrichia curtipendula'
$ws.Range("H12").Value = '(Hedw.) Brid.'
$ws.Range("Q12").Value = 658722.9088558007
$ws.Range("R12").Value = 6636991.191442309

# Row 13
$ws.Range("A13").Value = 96779993
$ws.Range("B13").Value = 43464
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 101735
$ws.Range("F13").Value = 'Jättesvampmal'
$ws.Range("G13").Value = 'Scardia boletella'
$ws.Range("H13").Value = '(Fabricius, 1794)'
$ws.Range("Q13").Value = 658739.4120713262
$ws.Range("R13").Value = 6636888.229354058

# Row 14
$ws.Range("A14").Value = 96780357
$ws.Range("B14").Value = 98520
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 222498
$ws.Range("F14").Value = 'Blåsippa'
$ws.Range("G14").Value = 'Hepatica nobilis'
$ws.Range("H14").Value = 'Schreb.'
$ws.Range("Q14").Value = 658747.5451754113
$ws.Range("R14").Value = 6637110.504147635

# Row 15
$ws.Range("A15").Value = 96780175
$ws.Range("B15").Value = 90005
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 1339
$ws.Range("F15").Value = 'Brandticka'
$ws.Range("G15").Value = 'Pycnoporellus fulgens'
$ws.Range("H15").Value = '(Fr.) Donk'
$ws.Range("Q15").Value = 658742.8409314866
$ws.Range("R15").Value = 6636937.694258579

# Activity column (M) updates
$ws.Range("M11").Value = 'äldre gnagspår'
$ws.Range("M13").Value = 'äldre gnagspår'
$ws.Range("M8").Value = ""
$ws.Range("M10").Value = ""
